$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2
$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("C4").Value = 6
$ws.Range("E4").Value = 19
$ws.Range("F4").Value = 2
$ws.Range("D5").Value = 16
$ws.Range("F5").Value = 9
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("F7").Value = 1
